$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking Price cells so they are not
# auto-converted to numbers (preserves exact text representation).
$textCells = @('D5', 'D6', 'D10', 'D13', 'D14', 'D20', 'D21', 'D22', 'D23', 'D24', 'D25', 'D27', 'D29', 'D32', 'D33', 'D36', 'D37', 'D38', 'D39', 'D41', 'D44', 'D46', 'D47', 'D48', 'D51')
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '62.735.11'
$ws.Range('E2').Value = '  +1.85%  '
$ws.Range('D3').Value = '2.938.62'
$ws.Range('E3').Value = '  -0.05%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '592.85'
$ws.Range('E5').Value = '  -1.02%  '
$ws.Range('D6').Value = '146.91'
$ws.Range('E6').Value = '  +1.04%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('E8').Value = '  +0.70%  '
$ws.Range('D9').Value = '2.939.21'
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').Value = '7.31'
$ws.Range('E10').Value = '  +4.57%  '
$ws.Range('E11').Value = '  +5.85%  '
$ws.Range('E12').Value = '  +0.27%  '
$ws.Range('D13').Value = '0.0000237'
$ws.Range('E13').Value = '  +5.08%  '
$ws.Range('D14').Value = '32.75'
$ws.Range('E14').Value = '  -2.90%  '
$ws.Range('E15').Value = '  -0.97%  '
$ws.Range('D16').Value = '3.424.59'
$ws.Range('E16').Value = '  +0.17%  '
$ws.Range('D17').Value = '62.665.69'
$ws.Range('E17').Value = '  +1.84%  '
$ws.Range('E18').Value = '  -0.50%  '
$ws.Range('D19').Value = '2.957.35'
$ws.Range('E19').Value = '  +0.80%  '
$ws.Range('D20').Value = '440.88'
$ws.Range('E20').Value = '  +1.50%  '
$ws.Range('D21').Value = '13.39'
$ws.Range('E21').Value = '  -0.88%  '
$ws.Range('D22').Value = '0.665'
$ws.Range('E22').Value = '  -2.15%  '
$ws.Range('D23').Value = '7.02'
$ws.Range('E23').Value = '  -1.42%  '
$ws.Range('D24').Value = '81.26'
$ws.Range('E24').Value = '  -0.76%  '
$ws.Range('D25').Value = '11.09'
$ws.Range('E25').Value = '  +1.45%  '
$ws.Range('E26').Value = '  -3.27%  '
$ws.Range('D27').Value = '11.72'
$ws.Range('E27').Value = '  -0.79%  '
$ws.Range('E28').Value = '  -0.05%  '
$ws.Range('D29').Value = '2.24'
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('E30').Value = '  +3.87%  '
$ws.Range('E31').Value = '  -0.60%  '
$ws.Range('D32').Value = '0.0000101'
$ws.Range('E32').Value = '  +14.89%  '
$ws.Range('D33').Value = '26.35'
$ws.Range('E33').Value = '  -1.46%  '
$ws.Range('E34').Value = '  -1.44%  '
$ws.Range('E35').Value = '  +0.30%  '
$ws.Range('D36').Value = '0.988'
$ws.Range('E36').Value = '  -2.36%  '
$ws.Range('D37').Value = '3.12'
$ws.Range('E37').Value = '  +3.74%  '
$ws.Range('D38').Value = '5.58'
$ws.Range('E38').Value = '  -1.38%  '
$ws.Range('D39').Value = '49.66'
$ws.Range('E39').Value = '  -0.59%  '
$ws.Range('E40').Value = '  +0.64%  '
$ws.Range('D41').Value = '8.47'
$ws.Range('E41').Value = '  -1.76%  '
$ws.Range('E42').Value = '  -5.37%  '
$ws.Range('E43').Value = '  -1.39%  '
$ws.Range('D44').Value = '39.24'
$ws.Range('E44').Value = '  -7.51%  '
$ws.Range('D45').Value = '2.700.25'
$ws.Range('E45').Value = '  -0.11%  '
$ws.Range('D46').Value = '134.81'
$ws.Range('E46').Value = '  +0.16%  '
$ws.Range('B47').Value = 'Bittensor'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D47').Value = '362.03'
$ws.Range('E47').Value = '  -0.88%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').Value = '0.0336'
$ws.Range('E48').Value = '  -3.42%  '
$ws.Range('E50').Value = '  -0.83%  '
$ws.Range('D51').Value = '22.78'
$ws.Range('E51').Value = '  -4.33%  '
